# Update NATMI ligand-receptor pair table (Col2a1-Tnfrsf10b) with values
# recomputed from the refreshed TPM expression matrix ("update scripts wuth new tpm").
# Ligand (Col2a1) stats per sending cluster: columns E-J
# Receptor (Tnfrsf10b) stats per target cluster: columns M-P
# Edge weights/specificities derived from the above: columns Q-T

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> ECs ---
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.1581976666666667
$ws.Range("H2").Value = 0.474593
$ws.Range("I2").Value = 0.1400666049254827
$ws.Range("J2").Value = 0.1400666049254826
$ws.Range("M2").Value = 14.10125566666667
$ws.Range("N2").Value = 42.303767
$ws.Range("O2").Value = 0.7585903740943118
$ws.Range("P2").Value = 0.7585903740943116
$ws.Range("Q2").Value = 2.230785743536778
$ws.Range("R2").Value = 20.077071691831
$ws.Range("S2").Value = 0.1062531782285421
$ws.Range("T2").Value = 0.106253178228542

# --- Row 3: ECs -> FAPs ---
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.1581976666666667
$ws.Range("H3").Value = 0.474593
$ws.Range("I3").Value = 0.1400666049254827
$ws.Range("J3").Value = 0.1400666049254826
$ws.Range("O3").Value = 0.1642962051354147
$ws.Range("P3").Value = 0.1642962051354147
$ws.Range("Q3").Value = 0.483145640453
$ws.Range("R3").Value = 4.348310764077
$ws.Range("S3").Value = 0.02301241165545819
$ws.Range("T3").Value = 0.02301241165545818

# --- Row 4: ECs -> MuSCs ---
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.1581976666666667
$ws.Range("H4").Value = 0.474593
$ws.Range("I4").Value = 0.1400666049254827
$ws.Range("J4").Value = 0.1400666049254826
$ws.Range("M4").Value = 1.366842
$ws.Range("N4").Value = 4.100526
$ws.Range("O4").Value = 0.07353055703818179
$ws.Range("P4").Value = 0.07353055703818176
$ws.Range("Q4").Value = 0.216231215102
$ws.Range("R4").Value = 1.946080935918
$ws.Range("S4").Value = 0.01029917548261768
$ws.Range("T4").Value = 0.01029917548261767

# --- Row 5: ECs -> Resolving-Mac ---
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.1581976666666667
$ws.Range("H5").Value = 0.474593
$ws.Range("I5").Value = 0.1400666049254827
$ws.Range("J5").Value = 0.1400666049254826
$ws.Range("M5").Value = 0.06660100000000001
$ws.Range("N5").Value = 0.199803
$ws.Range("O5").Value = 0.003582863732091891
$ws.Range("P5").Value = 0.00358286373209189
$ws.Range("Q5").Value = 0.01053612279766667
$ws.Range("R5").Value = 0.09482510517900002
$ws.Range("S5").Value = 0.0005018395588647554
$ws.Range("T5").Value = 0.0005018395588647551

# --- Row 6: FAPs -> ECs ---
$ws.Range("I6").Value = 0.7029419733214338
$ws.Range("J6").Value = 0.7029419733214337
$ws.Range("M6").Value = 14.10125566666667
$ws.Range("N6").Value = 42.303767
$ws.Range("O6").Value = 0.7585903740943118
$ws.Range("P6").Value = 0.7585903740943116
$ws.Range("Q6").Value = 11.195480417715
$ws.Range("R6").Value = 100.759323759435
$ws.Range("S6").Value = 0.5332450145085003
$ws.Range("T6").Value = 0.5332450145085

# --- Row 7: FAPs -> FAPs ---
$ws.Range("I7").Value = 0.7029419733214338
$ws.Range("J7").Value = 0.7029419733214337
$ws.Range("O7").Value = 0.1642962051354147
$ws.Range("P7").Value = 0.1642962051354147
$ws.Range("S7").Value = 0.1154906986471115
$ws.Range("T7").Value = 0.1154906986471115

# --- Row 8: FAPs -> MuSCs ---
$ws.Range("I8").Value = 0.7029419733214338
$ws.Range("J8").Value = 0.7029419733214337
$ws.Range("M8").Value = 1.366842
$ws.Range("N8").Value = 4.100526
$ws.Range("O8").Value = 0.07353055703818179
$ws.Range("P8").Value = 0.07353055703818176
$ws.Range("Q8").Value = 1.08518370327
$ws.Range("R8").Value = 9.766653329430001
$ws.Range("S8").Value = 0.05168771486384375
$ws.Range("T8").Value = 0.05168771486384372

# --- Row 9: FAPs -> Resolving-Mac ---
$ws.Range("I9").Value = 0.7029419733214338
$ws.Range("J9").Value = 0.7029419733214337
$ws.Range("M9").Value = 0.06660100000000001
$ws.Range("N9").Value = 0.199803
$ws.Range("O9").Value = 0.003582863732091891
$ws.Range("P9").Value = 0.00358286373209189
$ws.Range("Q9").Value = 0.052876864935
$ws.Range("R9").Value = 0.475891784415
$ws.Range("S9").Value = 0.002518545301978471
$ws.Range("T9").Value = 0.00251854530197847

# --- Row 10: MuSCs -> ECs ---
$ws.Range("G10").Value = 0.1199896666666667
$ws.Range("H10").Value = 0.359969
$ws.Range("I10").Value = 0.106237630366274
$ws.Range("J10").Value = 0.106237630366274
$ws.Range("M10").Value = 14.10125566666667
$ws.Range("N10").Value = 42.303767
$ws.Range("O10").Value = 0.7585903740943118
$ws.Range("P10").Value = 0.7585903740943116
$ws.Range("Q10").Value = 1.692004967024778
$ws.Range("R10").Value = 15.228044703223
$ws.Range("S10").Value = 0.080590843762445
$ws.Range("T10").Value = 0.08059084376244496

# --- Row 11: MuSCs -> FAPs ---
$ws.Range("G11").Value = 0.1199896666666667
$ws.Range("H11").Value = 0.359969
$ws.Range("I11").Value = 0.106237630366274
$ws.Range("J11").Value = 0.106237630366274
$ws.Range("O11").Value = 0.1642962051354147
$ws.Range("P11").Value = 0.1642962051354147
$ws.Range("Q11").Value = 0.366456001349
$ws.Range("R11").Value = 3.298104012141
$ws.Range("S11").Value = 0.01745443951175771
$ws.Range("T11").Value = 0.0174544395117577

# --- Row 12: MuSCs -> MuSCs ---
$ws.Range("G12").Value = 0.1199896666666667
$ws.Range("H12").Value = 0.359969
$ws.Range("I12").Value = 0.106237630366274
$ws.Range("J12").Value = 0.106237630366274
$ws.Range("M12").Value = 1.366842
$ws.Range("N12").Value = 4.100526
$ws.Range("O12").Value = 0.07353055703818179
$ws.Range("P12").Value = 0.07353055703818176
$ws.Range("Q12").Value = 0.164006915966
$ws.Range("R12").Value = 1.476062243694
$ws.Range("S12").Value = 0.007811712139248582
$ws.Range("T12").Value = 0.007811712139248577

# --- Row 13: MuSCs -> Resolving-Mac ---
$ws.Range("G13").Value = 0.1199896666666667
$ws.Range("H13").Value = 0.359969
$ws.Range("I13").Value = 0.106237630366274
$ws.Range("J13").Value = 0.106237630366274
$ws.Range("M13").Value = 0.06660100000000001
$ws.Range("N13").Value = 0.199803
$ws.Range("O13").Value = 0.003582863732091891
$ws.Range("P13").Value = 0.00358286373209189
$ws.Range("Q13").Value = 0.007991431789666667
$ws.Range("R13").Value = 0.071922886107
$ws.Range("S13").Value = 0.0003806349528227073
$ws.Range("T13").Value = 0.000380634952822707

# --- Row 14: Resolving-Mac -> ECs ---
$ws.Range("E14").Value = 2.0
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05732366666666666
$ws.Range("H14").Value = 0.171971
$ws.Range("I14").Value = 0.05075379138680971
$ws.Range("J14").Value = 0.05075379138680969
$ws.Range("M14").Value = 14.10125566666667
$ws.Range("N14").Value = 42.303767
$ws.Range("O14").Value = 0.7585903740943118
$ws.Range("P14").Value = 0.7585903740943116
$ws.Range("Q14").Value = 0.8083356794174443
$ws.Range("R14").Value = 7.275021114756999
$ws.Range("S14").Value = 0.03850133759482464
$ws.Range("T14").Value = 0.03850133759482462

# --- Row 15: Resolving-Mac -> FAPs ---
$ws.Range("E15").Value = 2.0
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05732366666666666
$ws.Range("H15").Value = 0.171971
$ws.Range("I15").Value = 0.05075379138680971
$ws.Range("J15").Value = 0.05075379138680969
$ws.Range("O15").Value = 0.1642962051354147
$ws.Range("P15").Value = 0.1642962051354147
$ws.Range("Q15").Value = 0.175070089391
$ws.Range("R15").Value = 1.575630804519
$ws.Range("S15").Value = 0.008338655321087334
$ws.Range("T15").Value = 0.008338655321087329

# --- Row 16: Resolving-Mac -> MuSCs ---
$ws.Range("E16").Value = 2.0
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05732366666666666
$ws.Range("H16").Value = 0.171971
$ws.Range("I16").Value = 0.05075379138680971
$ws.Range("J16").Value = 0.05075379138680969
$ws.Range("M16").Value = 1.366842
$ws.Range("N16").Value = 4.100526
$ws.Range("O16").Value = 0.07353055703818179
$ws.Range("P16").Value = 0.07353055703818176
$ws.Range("Q16").Value = 0.078352395194
$ws.Range("R16").Value = 0.705171556746
$ws.Range("S16").Value = 0.003731954552471791
$ws.Range("T16").Value = 0.003731954552471788

# --- Row 17: Resolving-Mac -> Resolving-Mac ---
$ws.Range("E17").Value = 2.0
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05732366666666666
$ws.Range("H17").Value = 0.171971
$ws.Range("I17").Value = 0.05075379138680971
$ws.Range("J17").Value = 0.05075379138680969
$ws.Range("M17").Value = 0.06660100000000001
$ws.Range("N17").Value = 0.199803
$ws.Range("O17").Value = 0.003582863732091891
$ws.Range("P17").Value = 0.00358286373209189
$ws.Range("Q17").Value = 0.003817813523666667
$ws.Range("R17").Value = 0.034360321713
$ws.Range("S17").Value = 0.0001818439184259583
$ws.Range("T17").Value = 0.0001818439184259582
